$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the "Template" sheet to "res.partner"
$ws.Name = "res.partner"

# --- Header row (row 1): switch from human-friendly labels to Odoo technical field names ---
$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "company_type"
$ws.Range("D1").Value = "parent_id"
$ws.Range("E1").Value = "type"
$ws.Range("F1").Value = "street"
$ws.Range("G1").Value = "street2"
$ws.Range("H1").Value = "city"
$ws.Range("I1").Value = "state_id/id"
$ws.Range("J1").Value = "zip"
$ws.Range("K1").Value = "country_id"
$ws.Range("L1").Value = "website"
$ws.Range("M1").Value = "phone"
$ws.Range("N1").Value = "mobile"
$ws.Range("O1").Value = "email"

# --- Row 2: XLSX Grocery Outlet ---
$ws.Range("A2").Value = "__import__.res_partner_1"
$ws.Range("I2").Value = "base.state_us_5"

# --- Row 3: XLSX Bob ---
$ws.Range("A3").Value = "__import__.res_partner_2"
$ws.Range("D3").Value = "XLSX Grocery Outlet"

# --- Row 4: XLSX Bill ---
$ws.Range("A4").Value = "__import__.res_partner_3"
$ws.Range("D4").Value = "XLSX Grocery Outlet"
$ws.Range("I4").Value = "base.state_us_5"

# --- Row 5: XLSX Warehouse ---
$ws.Range("A5").Value = "__import__.res_partner_4"
$ws.Range("D5").Value = "XLSX Grocery Outlet"
$ws.Range("I5").Value = "base.state_us_5"

# --- Row 6: XLSX Accounting dep. ---
$ws.Range("A6").Value = "__import__.res_partner_5"
$ws.Range("D6").Value = "XLSX Grocery Outlet"
$ws.Range("I6").Value = "base.state_us_5"

# --- Row 7: XLSX ASUSTeK ---
$ws.Range("A7").Value = "__import__.res_partner_6"
$ws.Range("I7").Value = "base.state_us_27"

# --- Row 8: XLSX Camptocamp ---
$ws.Range("A8").Value = "__import__.res_partner_7"

# --- Column widths: column A gets its own wider width to fit the longer
#     technical id strings (columns B/C keep their existing 14.13 width), and
#     column D is nudged by the tiny fractional amount produced alongside
#     the resize (20.19 -> 20.18). The COM layer quantizes ColumnWidth to a
#     pixel grid, so the inputs below are tuned to land on the closest
#     achievable stored width (33.17 / 20.17 vs. targets 33.2 / 20.18) ---
$ws.Columns.Item(1).ColumnWidth = 32.25
$ws.Columns.Item(4).ColumnWidth = 19.25

# --- Selection moves from B8 to I1 ---
$ws.Range("I1").Select()
